# Shift the GESS model forward by 5 days: every timestamp in column A (Timestamp)
# and every lookup key in column E (Lookup) moves from the 19-21 Feb 2026 window
# to the 24-26 Feb 2026 window, keeping the embedded Quarter counter intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 46077.99444444444
$ws.Cells.Item(2, 5).Value = "24.02.20261"
$ws.Cells.Item(3, 1).Value = 46078.00486111111
$ws.Cells.Item(3, 5).Value = "25.02.20262"
$ws.Cells.Item(4, 1).Value = 46078.01527777778
$ws.Cells.Item(4, 5).Value = "25.02.20263"
$ws.Cells.Item(5, 1).Value = 46078.02569444444
$ws.Cells.Item(5, 5).Value = "25.02.20264"
$ws.Cells.Item(6, 1).Value = 46078.03611111111
$ws.Cells.Item(6, 5).Value = "25.02.20265"
$ws.Cells.Item(7, 1).Value = 46078.04652777778
$ws.Cells.Item(7, 5).Value = "25.02.20266"
$ws.Cells.Item(8, 1).Value = 46078.05694444444
$ws.Cells.Item(8, 5).Value = "25.02.20267"
$ws.Cells.Item(9, 1).Value = 46078.06736111111
$ws.Cells.Item(9, 5).Value = "25.02.20268"
$ws.Cells.Item(10, 1).Value = 46078.07777777778
$ws.Cells.Item(10, 5).Value = "25.02.20269"
$ws.Cells.Item(11, 1).Value = 46078.08819444444
$ws.Cells.Item(11, 5).Value = "25.02.202610"
$ws.Cells.Item(12, 1).Value = 46078.09861111111
$ws.Cells.Item(12, 5).Value = "25.02.202611"
$ws.Cells.Item(13, 1).Value = 46078.10902777778
$ws.Cells.Item(13, 5).Value = "25.02.202612"
$ws.Cells.Item(14, 1).Value = 46078.11944444444
$ws.Cells.Item(14, 5).Value = "25.02.202613"
$ws.Cells.Item(15, 1).Value = 46078.12986111111
$ws.Cells.Item(15, 5).Value = "25.02.202614"
$ws.Cells.Item(16, 1).Value = 46078.14027777778
$ws.Cells.Item(16, 5).Value = "25.02.202615"
$ws.Cells.Item(17, 1).Value = 46078.15069444444
$ws.Cells.Item(17, 5).Value = "25.02.202616"
$ws.Cells.Item(18, 1).Value = 46078.16111111111
$ws.Cells.Item(18, 5).Value = "25.02.202617"
$ws.Cells.Item(19, 1).Value = 46078.17152777778
$ws.Cells.Item(19, 5).Value = "25.02.202618"
$ws.Cells.Item(20, 1).Value = 46078.18194444444
$ws.Cells.Item(20, 5).Value = "25.02.202619"
$ws.Cells.Item(21, 1).Value = 46078.19236111111
$ws.Cells.Item(21, 5).Value = "25.02.202620"
$ws.Cells.Item(22, 1).Value = 46078.20277777778
$ws.Cells.Item(22, 5).Value = "25.02.202621"
$ws.Cells.Item(23, 1).Value = 46078.21319444444
$ws.Cells.Item(23, 5).Value = "25.02.202622"
$ws.Cells.Item(24, 1).Value = 46078.22361111111
$ws.Cells.Item(24, 5).Value = "25.02.202623"
$ws.Cells.Item(25, 1).Value = 46078.23402777778
$ws.Cells.Item(25, 5).Value = "25.02.202624"
$ws.Cells.Item(26, 1).Value = 46078.24444444444
$ws.Cells.Item(26, 5).Value = "25.02.202625"
$ws.Cells.Item(27, 1).Value = 46078.25486111111
$ws.Cells.Item(27, 5).Value = "25.02.202626"
$ws.Cells.Item(28, 1).Value = 46078.26527777778
$ws.Cells.Item(28, 5).Value = "25.02.202627"
$ws.Cells.Item(29, 1).Value = 46078.27569444444
$ws.Cells.Item(29, 5).Value = "25.02.202628"
$ws.Cells.Item(30, 1).Value = 46078.28611111111
$ws.Cells.Item(30, 5).Value = "25.02.202629"
$ws.Cells.Item(31, 1).Value = 46078.29652777778
$ws.Cells.Item(31, 5).Value = "25.02.202630"
$ws.Cells.Item(32, 1).Value = 46078.30694444444
$ws.Cells.Item(32, 5).Value = "25.02.202631"
$ws.Cells.Item(33, 1).Value = 46078.31736111111
$ws.Cells.Item(33, 5).Value = "25.02.202632"
$ws.Cells.Item(34, 1).Value = 46078.32777777778
$ws.Cells.Item(34, 5).Value = "25.02.202633"
$ws.Cells.Item(35, 1).Value = 46078.33819444444
$ws.Cells.Item(35, 5).Value = "25.02.202634"
$ws.Cells.Item(36, 1).Value = 46078.34861111111
$ws.Cells.Item(36, 5).Value = "25.02.202635"
$ws.Cells.Item(37, 1).Value = 46078.35902777778
$ws.Cells.Item(37, 5).Value = "25.02.202636"
$ws.Cells.Item(38, 1).Value = 46078.36944444444
$ws.Cells.Item(38, 5).Value = "25.02.202637"
$ws.Cells.Item(39, 1).Value = 46078.37986111111
$ws.Cells.Item(39, 5).Value = "25.02.202638"
$ws.Cells.Item(40, 1).Value = 46078.39027777778
$ws.Cells.Item(40, 5).Value = "25.02.202639"
$ws.Cells.Item(41, 1).Value = 46078.40069444444
$ws.Cells.Item(41, 5).Value = "25.02.202640"
$ws.Cells.Item(42, 1).Value = 46078.41111111111
$ws.Cells.Item(42, 5).Value = "25.02.202641"
$ws.Cells.Item(43, 1).Value = 46078.42152777778
$ws.Cells.Item(43, 5).Value = "25.02.202642"
$ws.Cells.Item(44, 1).Value = 46078.43194444444
$ws.Cells.Item(44, 5).Value = "25.02.202643"
$ws.Cells.Item(45, 1).Value = 46078.44236111111
$ws.Cells.Item(45, 5).Value = "25.02.202644"
$ws.Cells.Item(46, 1).Value = 46078.45277777778
$ws.Cells.Item(46, 5).Value = "25.02.202645"
$ws.Cells.Item(47, 1).Value = 46078.46319444444
$ws.Cells.Item(47, 5).Value = "25.02.202646"
$ws.Cells.Item(48, 1).Value = 46078.47361111111
$ws.Cells.Item(48, 5).Value = "25.02.202647"
$ws.Cells.Item(49, 1).Value = 46078.48402777778
$ws.Cells.Item(49, 5).Value = "25.02.202648"
$ws.Cells.Item(50, 1).Value = 46078.49444444444
$ws.Cells.Item(50, 5).Value = "25.02.202649"
$ws.Cells.Item(51, 1).Value = 46078.50486111111
$ws.Cells.Item(51, 5).Value = "25.02.202650"
$ws.Cells.Item(52, 1).Value = 46078.51527777778
$ws.Cells.Item(52, 5).Value = "25.02.202651"
$ws.Cells.Item(53, 1).Value = 46078.52569444444
$ws.Cells.Item(53, 5).Value = "25.02.202652"
$ws.Cells.Item(54, 1).Value = 46078.53611111111
$ws.Cells.Item(54, 5).Value = "25.02.202653"
$ws.Cells.Item(55, 1).Value = 46078.54652777778
$ws.Cells.Item(55, 5).Value = "25.02.202654"
$ws.Cells.Item(56, 1).Value = 46078.55694444444
$ws.Cells.Item(56, 5).Value = "25.02.202655"
$ws.Cells.Item(57, 1).Value = 46078.56736111111
$ws.Cells.Item(57, 5).Value = "25.02.202656"
$ws.Cells.Item(58, 1).Value = 46078.57777777778
$ws.Cells.Item(58, 5).Value = "25.02.202657"
$ws.Cells.Item(59, 1).Value = 46078.58819444444
$ws.Cells.Item(59, 5).Value = "25.02.202658"
$ws.Cells.Item(60, 1).Value = 46078.59861111111
$ws.Cells.Item(60, 5).Value = "25.02.202659"
$ws.Cells.Item(61, 1).Value = 46078.60902777778
$ws.Cells.Item(61, 5).Value = "25.02.202660"
$ws.Cells.Item(62, 1).Value = 46078.61944444444
$ws.Cells.Item(62, 5).Value = "25.02.202661"
$ws.Cells.Item(63, 1).Value = 46078.62986111111
$ws.Cells.Item(63, 5).Value = "25.02.202662"
$ws.Cells.Item(64, 1).Value = 46078.64027777778
$ws.Cells.Item(64, 5).Value = "25.02.202663"
$ws.Cells.Item(65, 1).Value = 46078.65069444444
$ws.Cells.Item(65, 5).Value = "25.02.202664"
$ws.Cells.Item(66, 1).Value = 46078.66111111111
$ws.Cells.Item(66, 5).Value = "25.02.202665"
$ws.Cells.Item(67, 1).Value = 46078.67152777778
$ws.Cells.Item(67, 5).Value = "25.02.202666"
$ws.Cells.Item(68, 1).Value = 46078.68194444444
$ws.Cells.Item(68, 5).Value = "25.02.202667"
$ws.Cells.Item(69, 1).Value = 46078.69236111111
$ws.Cells.Item(69, 5).Value = "25.02.202668"
$ws.Cells.Item(70, 1).Value = 46078.70277777778
$ws.Cells.Item(70, 5).Value = "25.02.202669"
$ws.Cells.Item(71, 1).Value = 46078.71319444444
$ws.Cells.Item(71, 5).Value = "25.02.202670"
$ws.Cells.Item(72, 1).Value = 46078.72361111111
$ws.Cells.Item(72, 5).Value = "25.02.202671"
$ws.Cells.Item(73, 1).Value = 46078.73402777778
$ws.Cells.Item(73, 5).Value = "25.02.202672"
$ws.Cells.Item(74, 1).Value = 46078.74444444444
$ws.Cells.Item(74, 5).Value = "25.02.202673"
$ws.Cells.Item(75, 1).Value = 46078.75486111111
$ws.Cells.Item(75, 5).Value = "25.02.202674"
$ws.Cells.Item(76, 1).Value = 46078.76527777778
$ws.Cells.Item(76, 5).Value = "25.02.202675"
$ws.Cells.Item(77, 1).Value = 46078.77569444444
$ws.Cells.Item(77, 5).Value = "25.02.202676"
$ws.Cells.Item(78, 1).Value = 46078.78611111111
$ws.Cells.Item(78, 5).Value = "25.02.202677"
$ws.Cells.Item(79, 1).Value = 46078.79652777778
$ws.Cells.Item(79, 5).Value = "25.02.202678"
$ws.Cells.Item(80, 1).Value = 46078.80694444444
$ws.Cells.Item(80, 5).Value = "25.02.202679"
$ws.Cells.Item(81, 1).Value = 46078.81736111111
$ws.Cells.Item(81, 5).Value = "25.02.202680"
$ws.Cells.Item(82, 1).Value = 46078.82777777778
$ws.Cells.Item(82, 5).Value = "25.02.202681"
$ws.Cells.Item(83, 1).Value = 46078.83819444444
$ws.Cells.Item(83, 5).Value = "25.02.202682"
$ws.Cells.Item(84, 1).Value = 46078.84861111111
$ws.Cells.Item(84, 5).Value = "25.02.202683"
$ws.Cells.Item(85, 1).Value = 46078.85902777778
$ws.Cells.Item(85, 5).Value = "25.02.202684"
$ws.Cells.Item(86, 1).Value = 46078.86944444444
$ws.Cells.Item(86, 5).Value = "25.02.202685"
$ws.Cells.Item(87, 1).Value = 46078.87986111111
$ws.Cells.Item(87, 5).Value = "25.02.202686"
$ws.Cells.Item(88, 1).Value = 46078.89027777778
$ws.Cells.Item(88, 5).Value = "25.02.202687"
$ws.Cells.Item(89, 1).Value = 46078.90069444444
$ws.Cells.Item(89, 5).Value = "25.02.202688"
$ws.Cells.Item(90, 1).Value = 46078.91111111111
$ws.Cells.Item(90, 5).Value = "25.02.202689"
$ws.Cells.Item(91, 1).Value = 46078.92152777778
$ws.Cells.Item(91, 5).Value = "25.02.202690"
$ws.Cells.Item(92, 1).Value = 46078.93194444444
$ws.Cells.Item(92, 5).Value = "25.02.202691"
$ws.Cells.Item(93, 1).Value = 46078.94236111111
$ws.Cells.Item(93, 5).Value = "25.02.202692"
$ws.Cells.Item(94, 1).Value = 46078.95277777778
$ws.Cells.Item(94, 5).Value = "25.02.202693"
$ws.Cells.Item(95, 1).Value = 46078.96319444444
$ws.Cells.Item(95, 5).Value = "25.02.202694"
$ws.Cells.Item(96, 1).Value = 46078.97361111111
$ws.Cells.Item(96, 5).Value = "25.02.202695"
$ws.Cells.Item(97, 1).Value = 46078.98402777778
$ws.Cells.Item(97, 5).Value = "25.02.202696"
$ws.Cells.Item(98, 1).Value = 46078.99444444444
$ws.Cells.Item(98, 5).Value = "25.02.20261"
$ws.Cells.Item(99, 1).Value = 46078.99444444444
$ws.Cells.Item(99, 5).Value = "25.02.20262"
$ws.Cells.Item(100, 1).Value = 46079.00486111111
$ws.Cells.Item(100, 5).Value = "26.02.20263"
$ws.Cells.Item(101, 1).Value = 46079.01527777778
$ws.Cells.Item(101, 5).Value = "26.02.20264"
$ws.Cells.Item(102, 1).Value = 46079.02569444444
$ws.Cells.Item(102, 5).Value = "26.02.20265"
$ws.Cells.Item(103, 1).Value = 46079.03611111111
$ws.Cells.Item(103, 5).Value = "26.02.20266"
$ws.Cells.Item(104, 1).Value = 46079.04652777778
$ws.Cells.Item(104, 5).Value = "26.02.20267"
$ws.Cells.Item(105, 1).Value = 46079.05694444444
$ws.Cells.Item(105, 5).Value = "26.02.20268"
$ws.Cells.Item(106, 1).Value = 46079.06736111111
$ws.Cells.Item(106, 5).Value = "26.02.20269"
$ws.Cells.Item(107, 1).Value = 46079.07777777778
$ws.Cells.Item(107, 5).Value = "26.02.202610"
$ws.Cells.Item(108, 1).Value = 46079.08819444444
$ws.Cells.Item(108, 5).Value = "26.02.202611"
$ws.Cells.Item(109, 1).Value = 46079.09861111111
$ws.Cells.Item(109, 5).Value = "26.02.202612"
$ws.Cells.Item(110, 1).Value = 46079.10902777778
$ws.Cells.Item(110, 5).Value = "26.02.202613"
$ws.Cells.Item(111, 1).Value = 46079.11944444444
$ws.Cells.Item(111, 5).Value = "26.02.202614"
$ws.Cells.Item(112, 1).Value = 46079.12986111111
$ws.Cells.Item(112, 5).Value = "26.02.202615"
$ws.Cells.Item(113, 1).Value = 46079.14027777778
$ws.Cells.Item(113, 5).Value = "26.02.202616"
$ws.Cells.Item(114, 1).Value = 46079.15069444444
$ws.Cells.Item(114, 5).Value = "26.02.202617"
$ws.Cells.Item(115, 1).Value = 46079.16111111111
$ws.Cells.Item(115, 5).Value = "26.02.202618"
$ws.Cells.Item(116, 1).Value = 46079.17152777778
$ws.Cells.Item(116, 5).Value = "26.02.202619"
$ws.Cells.Item(117, 1).Value = 46079.18194444444
$ws.Cells.Item(117, 5).Value = "26.02.202620"
$ws.Cells.Item(118, 1).Value = 46079.19236111111
$ws.Cells.Item(118, 5).Value = "26.02.202621"
$ws.Cells.Item(119, 1).Value = 46079.20277777778
$ws.Cells.Item(119, 5).Value = "26.02.202622"
$ws.Cells.Item(120, 1).Value = 46079.21319444444
$ws.Cells.Item(120, 5).Value = "26.02.202623"
$ws.Cells.Item(121, 1).Value = 46079.22361111111
$ws.Cells.Item(121, 5).Value = "26.02.202624"
$ws.Cells.Item(122, 1).Value = 46079.23402777778
$ws.Cells.Item(122, 5).Value = "26.02.202625"
$ws.Cells.Item(123, 1).Value = 46079.24444444444
$ws.Cells.Item(123, 5).Value = "26.02.202626"
$ws.Cells.Item(124, 1).Value = 46079.25486111111
$ws.Cells.Item(124, 5).Value = "26.02.202627"
$ws.Cells.Item(125, 1).Value = 46079.26527777778
$ws.Cells.Item(125, 5).Value = "26.02.202628"
$ws.Cells.Item(126, 1).Value = 46079.27569444444
$ws.Cells.Item(126, 5).Value = "26.02.202629"
$ws.Cells.Item(127, 1).Value = 46079.28611111111
$ws.Cells.Item(127, 5).Value = "26.02.202630"
$ws.Cells.Item(128, 1).Value = 46079.29652777778
$ws.Cells.Item(128, 5).Value = "26.02.202631"
$ws.Cells.Item(129, 1).Value = 46079.30694444444
$ws.Cells.Item(129, 5).Value = "26.02.202632"
$ws.Cells.Item(130, 1).Value = 46079.31736111111
$ws.Cells.Item(130, 5).Value = "26.02.202633"
$ws.Cells.Item(131, 1).Value = 46079.32777777778
$ws.Cells.Item(131, 5).Value = "26.02.202634"
$ws.Cells.Item(132, 1).Value = 46079.33819444444
$ws.Cells.Item(132, 5).Value = "26.02.202635"
$ws.Cells.Item(133, 1).Value = 46079.34861111111
$ws.Cells.Item(133, 5).Value = "26.02.202636"
$ws.Cells.Item(134, 1).Value = 46079.35902777778
$ws.Cells.Item(134, 5).Value = "26.02.202637"
$ws.Cells.Item(135, 1).Value = 46079.36944444444
$ws.Cells.Item(135, 5).Value = "26.02.202638"
$ws.Cells.Item(136, 1).Value = 46079.37986111111
$ws.Cells.Item(136, 5).Value = "26.02.202639"
$ws.Cells.Item(137, 1).Value = 46079.39027777778
$ws.Cells.Item(137, 5).Value = "26.02.202640"
$ws.Cells.Item(138, 1).Value = 46079.40069444444
$ws.Cells.Item(138, 5).Value = "26.02.202641"
$ws.Cells.Item(139, 1).Value = 46079.41111111111
$ws.Cells.Item(139, 5).Value = "26.02.202642"
$ws.Cells.Item(140, 1).Value = 46079.42152777778
$ws.Cells.Item(140, 5).Value = "26.02.202643"
$ws.Cells.Item(141, 1).Value = 46079.43194444444
$ws.Cells.Item(141, 5).Value = "26.02.202644"
$ws.Cells.Item(142, 1).Value = 46079.44236111111
$ws.Cells.Item(142, 5).Value = "26.02.202645"
$ws.Cells.Item(143, 1).Value = 46079.45277777778
$ws.Cells.Item(143, 5).Value = "26.02.202646"
$ws.Cells.Item(144, 1).Value = 46079.46319444444
$ws.Cells.Item(144, 5).Value = "26.02.202647"
$ws.Cells.Item(145, 1).Value = 46079.47361111111
$ws.Cells.Item(145, 5).Value = "26.02.202648"
$ws.Cells.Item(146, 1).Value = 46079.48402777778
$ws.Cells.Item(146, 5).Value = "26.02.202649"
$ws.Cells.Item(147, 1).Value = 46079.49444444444
$ws.Cells.Item(147, 5).Value = "26.02.202650"
$ws.Cells.Item(148, 1).Value = 46079.50486111111
$ws.Cells.Item(148, 5).Value = "26.02.202651"
$ws.Cells.Item(149, 1).Value = 46079.51527777778
$ws.Cells.Item(149, 5).Value = "26.02.202652"
$ws.Cells.Item(150, 1).Value = 46079.52569444444
$ws.Cells.Item(150, 5).Value = "26.02.202653"
$ws.Cells.Item(151, 1).Value = 46079.53611111111
$ws.Cells.Item(151, 5).Value = "26.02.202654"
$ws.Cells.Item(152, 1).Value = 46079.54652777778
$ws.Cells.Item(152, 5).Value = "26.02.202655"
$ws.Cells.Item(153, 1).Value = 46079.55694444444
$ws.Cells.Item(153, 5).Value = "26.02.202656"
$ws.Cells.Item(154, 1).Value = 46079.56736111111
$ws.Cells.Item(154, 5).Value = "26.02.202657"
$ws.Cells.Item(155, 1).Value = 46079.57777777778
$ws.Cells.Item(155, 5).Value = "26.02.202658"
$ws.Cells.Item(156, 1).Value = 46079.58819444444
$ws.Cells.Item(156, 5).Value = "26.02.202659"
$ws.Cells.Item(157, 1).Value = 46079.59861111111
$ws.Cells.Item(157, 5).Value = "26.02.202660"
$ws.Cells.Item(158, 1).Value = 46079.60902777778
$ws.Cells.Item(158, 5).Value = "26.02.202661"
$ws.Cells.Item(159, 1).Value = 46079.61944444444
$ws.Cells.Item(159, 5).Value = "26.02.202662"
$ws.Cells.Item(160, 1).Value = 46079.62986111111
$ws.Cells.Item(160, 5).Value = "26.02.202663"
$ws.Cells.Item(161, 1).Value = 46079.64027777778
$ws.Cells.Item(161, 5).Value = "26.02.202664"
$ws.Cells.Item(162, 1).Value = 46079.65069444444
$ws.Cells.Item(162, 5).Value = "26.02.202665"
$ws.Cells.Item(163, 1).Value = 46079.66111111111
$ws.Cells.Item(163, 5).Value = "26.02.202666"
$ws.Cells.Item(164, 1).Value = 46079.67152777778
$ws.Cells.Item(164, 5).Value = "26.02.202667"
$ws.Cells.Item(165, 1).Value = 46079.68194444444
$ws.Cells.Item(165, 5).Value = "26.02.202668"
$ws.Cells.Item(166, 1).Value = 46079.69236111111
$ws.Cells.Item(166, 5).Value = "26.02.202669"
$ws.Cells.Item(167, 1).Value = 46079.70277777778
$ws.Cells.Item(167, 5).Value = "26.02.202670"
$ws.Cells.Item(168, 1).Value = 46079.71319444444
$ws.Cells.Item(168, 5).Value = "26.02.202671"
$ws.Cells.Item(169, 1).Value = 46079.72361111111
$ws.Cells.Item(169, 5).Value = "26.02.202672"
$ws.Cells.Item(170, 1).Value = 46079.73402777778
$ws.Cells.Item(170, 5).Value = "26.02.202673"
$ws.Cells.Item(171, 1).Value = 46079.74444444444
$ws.Cells.Item(171, 5).Value = "26.02.202674"
$ws.Cells.Item(172, 1).Value = 46079.75486111111
$ws.Cells.Item(172, 5).Value = "26.02.202675"
$ws.Cells.Item(173, 1).Value = 46079.76527777778
$ws.Cells.Item(173, 5).Value = "26.02.202676"
$ws.Cells.Item(174, 1).Value = 46079.77569444444
$ws.Cells.Item(174, 5).Value = "26.02.202677"
$ws.Cells.Item(175, 1).Value = 46079.78611111111
$ws.Cells.Item(175, 5).Value = "26.02.202678"
$ws.Cells.Item(176, 1).Value = 46079.79652777778
$ws.Cells.Item(176, 5).Value = "26.02.202679"
$ws.Cells.Item(177, 1).Value = 46079.80694444444
$ws.Cells.Item(177, 5).Value = "26.02.202680"
$ws.Cells.Item(178, 1).Value = 46079.81736111111
$ws.Cells.Item(178, 5).Value = "26.02.202681"
$ws.Cells.Item(179, 1).Value = 46079.82777777778
$ws.Cells.Item(179, 5).Value = "26.02.202682"
$ws.Cells.Item(180, 1).Value = 46079.83819444444
$ws.Cells.Item(180, 5).Value = "26.02.202683"
$ws.Cells.Item(181, 1).Value = 46079.84861111111
$ws.Cells.Item(181, 5).Value = "26.02.202684"
$ws.Cells.Item(182, 1).Value = 46079.85902777778
$ws.Cells.Item(182, 5).Value = "26.02.202685"
$ws.Cells.Item(183, 1).Value = 46079.86944444444
$ws.Cells.Item(183, 5).Value = "26.02.202686"
$ws.Cells.Item(184, 1).Value = 46079.87986111111
$ws.Cells.Item(184, 5).Value = "26.02.202687"
$ws.Cells.Item(185, 1).Value = 46079.89027777778
$ws.Cells.Item(185, 5).Value = "26.02.202688"
$ws.Cells.Item(186, 1).Value = 46079.90069444444
$ws.Cells.Item(186, 5).Value = "26.02.202689"
$ws.Cells.Item(187, 1).Value = 46079.91111111111
$ws.Cells.Item(187, 5).Value = "26.02.202690"
$ws.Cells.Item(188, 1).Value = 46079.92152777778
$ws.Cells.Item(188, 5).Value = "26.02.202691"
$ws.Cells.Item(189, 1).Value = 46079.93194444444
$ws.Cells.Item(189, 5).Value = "26.02.202692"
$ws.Cells.Item(190, 1).Value = 46079.94236111111
$ws.Cells.Item(190, 5).Value = "26.02.202693"
$ws.Cells.Item(191, 1).Value = 46079.95277777778
$ws.Cells.Item(191, 5).Value = "26.02.202694"
$ws.Cells.Item(192, 1).Value = 46079.96319444444
$ws.Cells.Item(192, 5).Value = "26.02.202695"
$ws.Cells.Item(193, 1).Value = 46079.97361111111
$ws.Cells.Item(193, 5).Value = "26.02.202696"
$ws.Cells.Item(194, 1).Value = 46079.98402777778
$ws.Cells.Item(194, 5).Value = "26.02.20261"
$ws.Cells.Item(195, 1).Value = 46079.99444444444
$ws.Cells.Item(195, 5).Value = "26.02.20262"

Write-Host "Shifted rows 2-195 (columns A and E) forward by 5 days."
